$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, matching the authoritative diff.
# Values must be written as literal text (not auto-converted to number/percent)
# to preserve the original inline-string cell semantics, so we briefly force a
# Text number format for the write and then restore each cells original style.
$updates = [ordered]@{
    "D2" = "307.74"
    "E2" = "-4.61%"
    "D3" = "39.95"
    "E3" = "-7.04%"
    "D4" = "5.133"
    "E4" = "-1.11%"
    "D5" = "0.07731"
    "E5" = "-5.78%"
    "D6" = "4.246"
    "E6" = "-1.80%"
    "D7" = "1.625"
    "E7" = "-11.12%"
    "D8" = "0.8811"
    "E8" = "-5.55%"
    "E9" = "-9.06%"
    "D10" = "0.1753"
    "E10" = "-6.27%"
    "D11" = "0.08971"
    "E11" = "-5.33%"
    "D12" = "0.04399"
    "E12" = "-4.99%"
    "E13" = "-0.26%"
    "D14" = "0.001253"
    "E14" = "-4.43%"
    "D15" = "0.005797"
    "E15" = "-2.04%"
    "E16" = "-0.29%"
    "E17" = "-4.82%"
    "E18" = "-1.48%"
    "D19" = "7.008"
    "E19" = "-5.35%"
    "D20" = "0.1339"
    "E20" = "-3.41%"
    "E21" = "14.38%"
    "D22" = "0.04162"
    "E22" = "0.06%"
    "D23" = "0.001200"
    "E23" = "-3.46%"
    "D24" = "0.004103"
    "E24" = "-5.03%"
    "E25" = "8.41%"
    "E26" = "0.14%"
    "D38" = "0.02362"
    "E38" = "-14.56%"
    "D39" = "0.05148"
    "E39" = "-7.87%"
    "D40" = "0.007949"
    "E40" = "-1.60%"
    "D41" = "0.1328"
    "E41" = "-4.93%"
    "D42" = "0.006378"
    "E42" = "-2.50%"
    "D43" = "0.001943"
    "E43" = "-6.17%"
    "D44" = "0.008493"
    "E44" = "13.24%"
    "D45" = "0.3318"
    "E45" = "-4.57%"
    "D46" = "0.00006519"
    "E46" = "-6.81%"
    "E47" = "0.07%"
    "E48" = "98.58%"
    "D49" = "0.005827"
    "E49" = "67.57%"
    "D50" = "0.00002101"
    "E50" = "0.07%"
    "D51" = "0.0002001"
    "E51" = "0.07%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = $origStyle
}
